# Commit: Add *.py version + write function increase rate
#
# Sheet2!J4 was a rate counter that increases with each run; bump it
# from 0 to 127, and leave the selection where the user last left off
# (J5) on Sheet2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Activate()
$ws.Cells.Item(4, 10).Value = 127
$ws.Range("J5").Select()
